# Updates the "CommunicationOnBehalfOf" StructureDefinition summary workbook
# from FHIR IG version 5.0.0 (published 2021-12-16) to 6.0.0 (published
# 2022-01-21), swaps the Contact metadata row for Publisher/Jurisdiction
# info, and refreshes the root Extension row's Short/Definition text on the
# Elements sheet to reuse the SD's own Title/Description.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Metadata" sheet (Property / Value table)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Publication date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# The sheet used to list the same "Contact" / "No display for ContactDetail"
# row twice in a row (rows 10 and 11). Delete the second copy -- a pure
# structural shift that slides the untouched rows below (old rows 12-21,
# "Description" through "Context") up into rows 11-20 without re-typing
# their values, so cells like "false" (Abstract) keep their original string
# type instead of being re-parsed as a boolean by fresh input.
$meta.Rows.Item(11).Delete()

# Fill in the real Publisher value, and turn the former first "Contact" row
# (still sitting at row 10) into the new Jurisdiction row.
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---------------------------------------------------------------------------
# "Elements" sheet: the root Extension row now carries the SD's own title/
# description instead of the generic "Extension" / "An Extension" text.
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Communication On-Behalf Of"
$elements.Range("L2").Value = "Attributed provider, location, or organization communication is on-behalf of"
